$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Rows("143:143").Insert()
$ws.Range("A143").Value = 7
$ws.Range("B143").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C143").Value = "Ñuble"
$ws.Range("D143").Value = 44572
$ws.Range("E143").Value = 16
$ws.Range("F143").Value = 100112023
$ws.Range("G143").Value = "Brócoli"
$ws.Range("H143").Value = "Sin especificar"
$ws.Range("I143").Value = "Primera"
$ws.Range("J143").Value = 400
$ws.Range("K143").Value = 600
$ws.Range("L143").Value = 650
$ws.Range("M143").Value = 625
$ws.Range("N143").Value = "`$/unidad"
$ws.Range("O143").Value = "Región del Maule"
$ws.Range("P143").Value = 625
$ws.Range("Q143").Value = 1
$ws.Range("R143").Value = "Hortaliza"
